$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Insert a brand-new "2022-Q1" sheet right before "总计", after
#    "2021-Q4" -- same layout as the other quarterly fund-holding
#    sheets (2021-Q2 / 2021-Q3 / 2021-Q4).
# ---------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$newQ = $wb.Worksheets.Add($null, $q4)
$newQ.Name = "2022-Q1"

# Header row (bold, bordered, centered like the other quarter sheets)
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$cols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $newQ.Range($cols[$i] + "1")
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Row 2
$newQ.Range("A2").Value = 0
$newQ.Range("A2").Font.Bold = $true
$newQ.Range("A2").HorizontalAlignment = -4108
$newQ.Range("A2").VerticalAlignment = -4160
$newQ.Range("A2").Borders.LineStyle = 1

$newQ.Range("B2").NumberFormat = "@"
$newQ.Range("B2").Value = "008763"
$newQ.Range("C2").Value = "天弘越南市场股票（QDII）A"
$newQ.Range("D2").NumberFormat = "@"
$newQ.Range("D2").Value = "37.53"
$newQ.Range("E2").NumberFormat = "@"
$newQ.Range("E2").Value = "92.10"
$newQ.Range("F2").NumberFormat = "@"
$newQ.Range("F2").Value = "5.24"
$newQ.Range("G2").NumberFormat = "@"
$newQ.Range("G2").Value = "1.9666"
$newQ.Range("H2").Value = 6

# Row 3
$newQ.Range("A3").Value = 1
$newQ.Range("A3").Font.Bold = $true
$newQ.Range("A3").HorizontalAlignment = -4108
$newQ.Range("A3").VerticalAlignment = -4160
$newQ.Range("A3").Borders.LineStyle = 1

$newQ.Range("B3").NumberFormat = "@"
$newQ.Range("B3").Value = "008764"
$newQ.Range("C3").Value = "天弘越南市场股票（QDII）C"
$newQ.Range("D3").NumberFormat = "@"
$newQ.Range("D3").Value = "14.26"
$newQ.Range("E3").NumberFormat = "@"
$newQ.Range("E3").Value = "92.10"
$newQ.Range("F3").NumberFormat = "@"
$newQ.Range("F3").Value = "5.24"
$newQ.Range("G3").NumberFormat = "@"
$newQ.Range("G3").Value = "0.7472"
$newQ.Range("H3").Value = 6

# ---------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: add a new first data row for
#    2022-Q1 and bump the existing rows down.
# ---------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("A2").Font.Bold = $true
$total.Range("A2").HorizontalAlignment = -4108
$total.Range("A2").VerticalAlignment = -4160
$total.Range("A2").Borders.LineStyle = 1

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 2.71

# Renumber the rows that shifted down (0,1,2 -> 1,2,3)
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
